$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A").Delete()
